$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 321
$ws.Range("B2").Value = 654

$ws.Range("B2").Select()
